# Update financial plan for 6% amortizing loan
# Loan repayment doubles from 483.32 to 966.64/month, which flows through
# Total Expenses, Net Cash Flow and Ending Cash on the "Financial Plan"
# sheet, and into the summary metrics on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financial Plan")
$summary = $wb.Worksheets.Item("Summary")

$row13 = @(6799.97,6799.97,6799.97,6799.97,6799.97,6799.97,6799.97,6799.97,6799.97,6799.97,6799.97,6799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97,9799.97)
$row14 = @(-6150.07,-5630.15,-5070.24,-4550.32,-4030.4,-3485.48,-2950.57,-2405.65,-1885.73,-1350.82,-805.9,-285.98,-2341.61,-771.74,758.14,2313.02,3857.89,5412.77,6957.64,8512.52,10057.39,11612.27,13142.15,14712.02,16465.85,19430.69,22410.52,25375.36,28355.19,31335.02,34324.86,37304.69,40269.53,43249.36,46214.2,49194.03)
$row15 = @(43849.93,38219.78,33149.54,28599.22,24568.82,21083.34,18132.77,15727.12,13841.39,12490.57,11684.67,11398.69,9057.08,8285.34,9043.48,11356.5,15214.39,20627.16,27584.8,36097.32,46154.71,57766.98,70909.13,85621.15,102087.0,121517.69,143928.21,169303.57,197658.76,228993.78,263318.64,300623.33,340892.86,384142.22,430356.42,479550.45)

$loanRepayment = 966.64

for ($i = 0; $i -lt 36; $i++) {
    $col = $i + 2
    $ws.Cells.Item(16, $col).Value = $loanRepayment
    $ws.Cells.Item(13, $col).Value = $row13[$i]
    $ws.Cells.Item(14, $col).Value = $row14[$i]
    $ws.Cells.Item(15, $col).Value = $row15[$i]
}

# Summary sheet roll-up figures
$summary.Cells.Item(2, 2).Value = 966.64
$summary.Cells.Item(3, 2).Value = 34799.04
$summary.Cells.Item(4, 2).Value = 479550.45
$summary.Cells.Item(6, 2).Value = 8285.34

# Rows 5 (Opening Funding Balance) and 17 (Reserves) only apply to month 1;
# drop the stray blank numeric cells that trailed across the rest of the
# 36-month grid.
$ws.Range("C5:AK5").ClearContents()
$ws.Range("C17:AK17").ClearContents()
